$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New task rows 56-64 (Task_ID, Task name, and a single "P" marker column per row)
$ws.Range("A56").Value = 4.1
$ws.Range("B56").Value = "Complete Assessment Test"
$ws.Range("C56").Value = "P"

$ws.Range("A57").Value = 4.2
$ws.Range("B57").Value = "Added Encryption Service for Delicate Information"
$ws.Range("C57").Value = "P"

$ws.Range("A58").Value = 4.3
$ws.Range("B58").Value = "Manage Test Assessments"
$ws.Range("E58").Value = "P"

$ws.Range("A59").Value = 4.4
$ws.Range("B59").Value = "Manage Findings and Outcomes"
$ws.Range("E59").Value = "P"

$ws.Range("A60").Value = 4.5
$ws.Range("B60").Value = "Assign Patient Plan"
$ws.Range("D60").Value = "P"

$ws.Range("A61").Value = 4.6
$ws.Range("B61").Value = "Complete Exercises"
$ws.Range("D61").Value = "P"

$ws.Range("A62").Value = 4.7
$ws.Range("B62").Value = "Finished the Scheduler"
$ws.Range("G62").Value = "P"

$ws.Range("A63").Value = 4.8
$ws.Range("B63").Value = "Landing Pages - Home, About, etc.."
$ws.Range("F63").Value = "P"

$ws.Range("A64").Value = 4.9
$ws.Range("B64").Value = "Landing Pages - Physiotherapist, Admin"
$ws.Range("H64").Value = "P"

# Update the view state to match where the author ended up working
$ws.Range("A49").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H64").Select() | Out-Null
